# Project "Sample Project" save: update the 4th rule's name (cell B11 on
# the "Rules" sheet) from "R40" to "1".
#
# A plain $cell.Value = "1" assignment would be auto-typed as a number by
# Excel (since "1" parses as numeric), which would change the cell's
# storage type away from a shared-string text cell. To keep the cell a
# genuine text value (matching the original "R40" text cell), we write it
# as a formula that evaluates to the text string "1" and then convert
# that formula to a literal value in place via copy / paste-special
# (values only). This preserves the cell's existing style/number format
# while storing the result as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
